$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.862.87"
$ws.Range("E2").Value = "  -1.65%  "
$ws.Range("D3").Value = "1.829.12"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").Value = "'311.23"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").Value = "'1.007"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("D7").Value = "'0.4576"
$ws.Range("E7").Value = "  -0.80%  "
$ws.Range("D8").Value = "'0.3673"
$ws.Range("E8").Value = "  -0.96%  "
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("D10").Value = "'0.8735"
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("D11").Value = "'0.07808"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "'19.58"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D13").Value = "1.869.30"
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").Value = "'6.379"
$ws.Range("E15").Value = "  -2.59%  "
$ws.Range("D16").Value = "'87.16"
$ws.Range("E16").Value = "  -4.96%  "
$ws.Range("D17").Value = "'1.008"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").Value = "'0.000008718"
$ws.Range("E18").Value = "  -3.42%  "
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").Value = "26.898.26"
$ws.Range("E20").Value = "  -1.56%  "
$ws.Range("D21").Value = "'14.48"
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("D22").Value = "'4.989"
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("D23").Value = "'10.46"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").Value = "'2.000"
$ws.Range("E24").Value = "  +3.76%  "
$ws.Range("D25").Value = "'151.58"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").Value = "'18.21"
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("D27").Value = "'1.989"
$ws.Range("E27").Value = "  -3.98%  "
$ws.Range("D28").Value = "'113.82"
$ws.Range("E28").Value = "  -1.89%  "
$ws.Range("D29").Value = "'4.925"
$ws.Range("E29").Value = "  -3.46%  "
$ws.Range("E30").Value = "  -0.52%  "
$ws.Range("D31").Value = "'3.092"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("D32").Value = "'0.7406"
$ws.Range("E32").Value = "  -4.76%  "
$ws.Range("D33").Value = "'4.483"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("E34").Value = "  -3.63%  "
$ws.Range("D35").Value = "'2.510"
$ws.Range("E35").Value = "  -6.09%  "
$ws.Range("D36").Value = "'1.083"
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("D38").Value = "'0.05123"
$ws.Range("E38").Value = "  -2.12%  "
$ws.Range("D39").Value = "'2.907"
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("D40").Value = "'6.893"
$ws.Range("D41").Value = "'0.4967"
$ws.Range("E41").Value = "  -3.46%  "
$ws.Range("D42").Value = "'0.1592"
$ws.Range("E42").Value = "  -2.47%  "
$ws.Range("D43").Value = "'8.269"
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("D44").Value = "'0.4667"
$ws.Range("E44").Value = "  -3.05%  "
$ws.Range("D45").Value = "'1.008"
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'103.08"
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.07"
$ws.Range("E47").Value = "  -2.39%  "
$ws.Range("D48").Value = "'1.607"
$ws.Range("E48").Value = "  -2.58%  "
$ws.Range("D49").Value = "'0.06067"
$ws.Range("E49").Value = "  -2.42%  "
$ws.Range("D50").Value = "'64.67"
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("D51").Value = "'36.63"
$ws.Range("E51").Value = "  -0.51%  "
